$d = $word.ActiveDocument

# Locate the list-item paragraph that needs to be split in two.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Added ability to stack items*") {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    # Replace the paragraph's content with two ListParagraph items:
    #  1) "Added ability to stack items" (same text, but "items" wrapped
    #     in a grammar-check proofErr pair, split across two runs)
    #  2) a new sibling bullet "Added ability to use items" with the same
    #     style/numbering, "items" likewise wrapped in a proofErr pair.
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' +
          '</w:pPr>' +
          '<w:r><w:t xml:space="preserve">Added ability to stack </w:t></w:r>' +
          '<w:proofErr w:type="gramStart"/>' +
          '<w:r><w:t>items</w:t></w:r>' +
          '<w:proofErr w:type="gramEnd"/>' +
        '</w:p>' +
        '<w:p>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' +
          '</w:pPr>' +
          '<w:r><w:t xml:space="preserve">Added ability to use </w:t></w:r>' +
          '<w:proofErr w:type="gramStart"/>' +
          '<w:r><w:t>items</w:t></w:r>' +
          '<w:proofErr w:type="gramEnd"/>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $targetPara.Range.InsertXML($xmlFrag) | Out-Null
}
